$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a new delivery-list entry added below the existing header (row 1)
# and data (row 2) rows. Copy row 2's formatting down to row 3 first so the
# new row picks up the same style (s="2": bordered / shaded / text format),
# then fill in the actual values column by column.
$ws.Range("A2:AN2").Copy()
$ws.Range("A3:AN3").PasteSpecial(-4122)

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "5739885993"
$ws.Range("C3").Value = "32100020883943"
$ws.Range("D3").Value = "우체국"
$ws.Range("E3").Value = "6094321250183"
$ws.Range("F3").Value = "분리배송불가"
$ws.Range("H3").Value = "2023-12-04"
$ws.Range("I3").Value = "2023-12-08 15:39:00"
$ws.Range("J3").Value = "2023-11-25 11:10:35"
$ws.Range("K3").Value = "베이비뵨 바운서 밸런스 소프트 메쉬 스카이 블루"
$ws.Range("L3").Value = "단일상품"
$ws.Range("M3").Value = "베이비뵨 바운서 밸런스 소프트 메쉬 스카이 블루"
$ws.Range("N3").Value = "6720117486"
$ws.Range("O3").Value = "82844419727"
$ws.Range("P3").Value = "베이비뵨 바운서 밸런스 소프트 메쉬 스카이 블루,단일상품"
$ws.Range("Q3").Value = "8434555784"
$ws.Range("S3").Value = "240000"
$ws.Range("T3").Value = "유료"
$ws.Range("U3").Value = "30000"
$ws.Range("V3").Value = "0"
$ws.Range("W3").Value = "1"
$ws.Range("X3").Value = "240000"
$ws.Range("Y3").Value = "한규택"
$ws.Range("Z3").Value = "010-****-****"
$ws.Range("AA3").Value = "Han Seungho"
$ws.Range("AB3").Value = "010-****-****"
$ws.Range("AC3").Value = "13596"
$ws.Range("AE3").Value = "문 앞"
$ws.Range("AH3").Value = "2023-12-09 19:43:00"
$ws.Range("AI3").Value = "2023-12-16 20:10:20"
$ws.Range("AJ3").Value = "P210023276057"
$ws.Range("AK3").Value = "010-5253-6534"
$ws.Range("AM3").Value = "아이폰앱"
$ws.Range("AN3").Value = "판매자 배송"

# Columns G, R, AD, AF, AG, AL stay blank for this row (same as the source).

# Move/scroll the active selection onto the newly-added row, matching the
# workbook's saved view state after the edit.
$ws.Range("A3").Select()
